$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @(
    @(1.02, 1.031045999123778, 1.044590174670016, 1.030425032668677, 1.049967282177634, 1.036903361857703, 1.036183887986969, 1.047360797067923, 1.033236082118503, 1.052722868187691, 1.037655387959864),
    @(1.02, 1.033948427943214, 1.04563755085187, 1.032960537080903, 1.051647025989925, 1.037221769682328, 1.038718095448415, 1.048218826217017, 1.035575306365899, 1.054212728220169, 1.040193194286583),
    @(1.02, 1.035814127633758, 1.046309554273485, 1.034590023710654, 1.052725461778131, 1.037423449315544, 1.040345875460647, 1.048767723731452, 1.037077559331293, 1.055167852552223, 1.041823285933141),
    @(1.02, 1.036595589328364, 1.04659071382624, 1.035272460125481, 1.05317684430909, 1.037507201735159, 1.041027388233157, 1.04899698873607, 1.037706448858958, 1.055567286990513, 1.042505766531959),
    @(1.02, 1.036726633745503, 1.046637843076433, 1.035386893778817, 1.053252517652323, 1.037521203749211, 1.041141654754523, 1.049035396387886, 1.037811888262569, 1.055634231641509, 1.042620195324891),
    @(1.02, 1.035824580765111, 1.046313316423787, 1.034599152581291, 1.05273150095234, 1.037424572470401, 1.040354992797726, 1.048770793018444, 1.037085972937692, 1.055173198019913, 1.041832416217885),
    @(1.02, 1.032029498825999, 1.044945335764924, 1.031284270771936, 1.050536735347097, 1.037011875206985, 1.037042868285639, 1.047652088355946, 1.034029031385951, 1.053228240153428, 1.038515588109084),
    @(1.02, 1.025243621770316, 1.042490129318935, 1.025354367104152, 1.046602748667842, 1.036250950363943, 1.031111111496399, 1.045631679804133, 1.028552067917115, 1.049731122466159, 1.032575407544822),
    @(1.02, 1.020648163103027, 1.040822087416287, 1.021336887638001, 1.043932998272698, 1.035720474040188, 1.027087784850466, 1.044250545250755, 1.024835707062346, 1.0473504905795, 1.02854636731357),
    @(1.02, 1.018640091294401, 1.040092117618291, 1.019580995224703, 1.042765256370827, 1.035485147149451, 1.025328243663596, 1.043644117010873, 1.023210053418285, 1.046307459861024, 1.026784327376392),
    @(1.02, 1.01789136159147, 1.039819793224421, 1.018926238077988, 1.042329696276946, 1.035396879850566, 1.024671961341926, 1.043417578401668, 1.022603654999204, 1.045918153191685, 1.026127113058548),
    @(1.02, 1.018052097031043, 1.0398782615855, 1.019066802180078, 1.042423208025924, 1.035415852420992, 1.024812860319516, 1.043466230206722, 1.022733846720946, 1.046001746551081, 1.026268212128848),
    @(1.02, 1.018578259677867, 1.040069631439361, 1.019526925162714, 1.04272929001625, 1.035477868482524, 1.025274051007406, 1.043625417596701, 1.023159981039895, 1.046275318241741, 1.026730057760414),
    @(1.02, 1.018902065951963, 1.040187383459859, 1.0198100826768, 1.042917636106889, 1.035515964806616, 1.025557843799298, 1.04372332730803, 1.023422194798715, 1.046443624470272, 1.027014253570622),
    @(1.02, 1.020781039540906, 1.040870368850718, 1.021453069271949, 1.044010246448041, 1.035735972422657, 1.027204184935854, 1.04429061314322, 1.024943242443585, 1.047419452170451, 1.02866293270043),
    @(1.02, 1.021954720744578, 1.041296709677662, 1.022479243938723, 1.04469244046172, 1.035872463187214, 1.028232163363054, 1.044644194192983, 1.025892891602855, 1.048028264740643, 1.029692370974907),
    @(1.02, 1.022637559705874, 1.041544646166331, 1.023076228151638, 1.045089224934647, 1.035951533764591, 1.028830092056421, 1.044849624387153, 1.026445225238629, 1.048382200694175, 1.030291148795602),
    @(1.02, 1.02287009660727, 1.041629061171077, 1.023279521452811, 1.045224328321619, 1.035978403153523, 1.02903368920752, 1.044919534498612, 1.026633291247787, 1.048502685915828, 1.030495035078008),
    @(1.02, 1.021828977574055, 1.041251044153745, 1.022369307756046, 1.044619364489193, 1.035857875178264, 1.028122044800477, 1.04460634200285, 1.02579116735558, 1.047963066683151, 1.029582096031341),
    @(1.02, 1.018423397230572, 1.040013310594189, 1.019391501295162, 1.042639206818546, 1.035459630029366, 1.02513831743305, 1.043578576523214, 1.023034566440561, 1.046194810428359, 1.026594131428814),
    @(1.02, 1.016265673174338, 1.039228252028236, 1.017504490037623, 1.041383713607479, 1.035204276848704, 1.023246599647251, 1.042924937413424, 1.021286535467739, 1.045072144926918, 1.024699727186776),
    @(1.02, 1.01741112406549, 1.039645083852718, 1.01850625958157, 1.042050284851037, 1.035340118378173, 1.024250958171845, 1.043272157501828, 1.022214636944858, 1.045668339183625, 1.025705512015676),
    @(1.02, 1.021885800899967, 1.041271680726666, 1.022418988018976, 1.044652387858065, 1.035864468546404, 1.028171807805461, 1.044623448272777, 1.025837137043552, 1.047992530504655, 1.029631929705499),
    @(1.02, 1.027010142794284, 1.043130273366778, 1.026898362100104, 1.047627901722717, 1.03645171530541, 1.032656390930536, 1.046159942489949, 1.029979132024453, 1.050643717407844, 1.034122881453)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = $cols[$j] + $row
        $ws.Range($addr).Value = $rowVals[$j]
    }
}

Write-Output "done"